# Edit Arbeitszeiten.xlsx:
#  - D49: 1 -> 2
#  - New row 50: Date 43404 (31.10.2018), "David", new "Was" entry, 1 hour
#  - Update sheet view (scroll position / selection) to match author's final state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D49 hours changed from 1 to 2
$ws.Range("D49").Value = 2

# Copy A49's date formatting down to A50 before writing the new date value
$ws.Range("A49").Copy($ws.Range("A50"))
$ws.Range("A50").Value = 43404

# New row of time-tracking data
$ws.Range("B50").Value = "David"
$ws.Range("C50").Value = "Added Foreign Key to Database + Auto Increments + Backup.sql generated"
$ws.Range("D50").Value = 1

# Update the saved view/selection state
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("C51").Select()
